# Adds a new "2020-04-24"-aligned forecast column (W) and a new forecast-origin
# row (35, dated 2020-05-08) to both the "cases" and "deaths" sheets, matching
# the staircase layout already used throughout the workbook.

function Set-TextValue($cell, $text) {
    # Force Excel to store the value as a plain text string (matching the
    # existing t="s" shared-string cells) instead of auto-converting
    # ISO-looking dates into a date serial number.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Update-ForecastSheet($ws, $newColValues, $b21Value) {
    # --- New column W (23) ---
    # Header: continues the existing date sequence (column header index
    # mirrors row-label index - 1), reusing the existing shared string.
    Set-TextValue $ws.Cells.Item(1, 23) "2020-04-24"

    # Rows 2-21: column W stays blank for these forecast-origin rows, but the
    # cell must still materialize as an explicit empty cell like its siblings.
    $blankTop = $ws.Range($ws.Cells.Item(2, 23), $ws.Cells.Item(21, 23))
    $blankTop.Font.Bold = $false

    # Rows 22-35: the new forecast values for column W.
    foreach ($r in ($newColValues.Keys | Sort-Object)) {
        $ws.Cells.Item($r, 23).Value = $newColValues[$r]
    }

    # Row 21, column B previously had no observation; it now does.
    $ws.Cells.Item(21, 2).Value = $b21Value

    # --- New row 35 ---
    Set-TextValue $ws.Cells.Item(35, 1) "2020-05-08"

    # Columns B-V (2-22) stay blank on the new row, materialized like the
    # rest of the table.
    $blankRow = $ws.Range($ws.Cells.Item(35, 2), $ws.Cells.Item(35, 22))
    $blankRow.Font.Bold = $false

    # Column W on row 35 holds a value too, already written above as part of
    # the $newColValues loop (row 35 is included in that table).
}

$wb = $excel.ActiveWorkbook

$casesNewCol = @{
    22 = 56328
    23 = 59703
    24 = 64092
    25 = 68172
    26 = 72916
    27 = 79219
    28 = 84990
    29 = 92133
    30 = 97435
    31 = 106011
    32 = 113078
    33 = 119679
    34 = 126092
    35 = 131844
}

$deathsNewCol = @{
    22 = 3931
    23 = 4243
    24 = 4655
    25 = 5048
    26 = 5513
    27 = 6140
    28 = 6733
    29 = 7478
    30 = 8053
    31 = 8982
    32 = 9774
    33 = 10531
    34 = 11280
    35 = 11965
}

$wsCases = $wb.Worksheets.Item("cases")
Update-ForecastSheet $wsCases $casesNewCol 52995

$wsDeaths = $wb.Worksheets.Item("deaths")
Update-ForecastSheet $wsDeaths $deathsNewCol 3670
